$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB mirrors column BA for all existing rows (1-83),
# copying both values and formatting (date format on row 1).
$ws.Range("BA1:BA83").Copy($ws.Range("BB1:BB83"))

# BB1 gets its own new date value (next quarter after BA1).
$ws.Range("BB1").Value = 45986

# New row 84: next date in column A, with a new EQUIPMENT eval value in BB.
$ws.Range("A83").Copy($ws.Range("A84"))
$ws.Range("A84").Value = 45884

$ws.Range("BB84").Value = 0.8783323788356512
